# Reorganization: the deck used to start with a standalone "Employe" class
# code-demo slide (sldId 900) followed by the Fonctions/Modules/Classes
# comparison slide (sldId 899). After the edit only the comparison slide
# remains, so we delete the old first slide.
$p = $ppt.ActivePresentation
$p.Slides.Item(1).Delete()
